# BP-359 Bank excel statements upload
#
# The bank-code column (B) stops being numeric: the values are re-entered
# as text ("199999", "288888", "388888", "488888", "588888") so they round
# trip through the sheet as shared strings. Columns A and C pick up the
# same text number format, and the transaction-amount column (E) switches
# to a fixed 2-decimal numeric format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and C (rows 2-6): apply the text ("@") number format used by
# the re-entered bank-code column so all three share the same style.
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("C2:C6").NumberFormat = "@"

# Column B (rows 2-6): re-enter the bank codes as text so Excel stores
# them as shared strings (with the "@" text format) instead of numbers.
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("B2").Value = "199999"
$ws.Range("B3").Value = "288888"
$ws.Range("B4").Value = "388888"
$ws.Range("B5").Value = "488888"
$ws.Range("B6").Value = "588888"

# Column E (rows 2-6): format the transaction amount with 2 decimals.
$ws.Range("E2:E6").NumberFormat = "0.00"

# Move the active selection from C6 to C3.
$ws.Range("C3").Select() | Out-Null

# Record an explicit (portrait) page setup for the sheet.
$ws.PageSetup.Orientation = 1
